# Added JB (username "bothma") on computer "phy-tglab11" to the folders sheet.
# This populates a new column N (one column to the right of the existing M)
# with the computer's settings, mirroring the layout of the other computer
# columns (B..M), where each row holds the value for the field named in
# column A of that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("N1").Value  = "phy-tglab11"
$ws.Range("N2").Value  = "bothma"
$ws.Range("N3").Value  = "Z:\LivemRNA\RawData"
$ws.Range("N4").Value  = "Z:\LivemRNA\FISHAnalysisData"
$ws.Range("N5").Value  = "C:\Users\bothma\Dropbox\LivemRNADatabase"
$ws.Range("N6").Value  = "C:\Users\bothma\Dropbox\LivemRNAData"
$ws.Range("N8").Value  = "C:\Users\bothma\Dropbox\MS2Pausing"
$ws.Range("N9").Value  = "Z:\LivemRNA\mRNADynamics"
$ws.Range("N10").Value = "Z:\FISHDrosophila\Analysis\schnitzcells"

# Reflect the view state shift (the split moved one column to the right, to
# just before the newly-added column, and the right pane's active cell is
# now the last cell that was just entered).
$excel.ActiveWindow.SplitColumn = 11
[void]$ws.Range("N10").Select()
